# "Encrypted messages in chat" - add day 17 number + a new timeline row
# (day 18, 20/3/2024, 40 hours, "Encoded chat messages"), then move the
# selection down to where the new data was entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 was missing its "Day" number (column A) - fill it in.
$ws.Range("A20").Value = 17
$ws.Range("A20").HorizontalAlignment = -4108
$ws.Range("A20").VerticalAlignment = -4108

# New row 21: Day 18, 20/3/2024, 40 hours, "Encoded chat messages".
$ws.Range("A21:D21").HorizontalAlignment = -4108
$ws.Range("A21:D21").VerticalAlignment = -4108

$ws.Range("A21").Value = 18
$ws.Range("B21").Value = "20/3/2024"
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = "Encoded chat messages"

# Scroll/select to reflect where the user ended up after the edit.
[void]$ws.Range("A4").Select()
[void]$ws.Range("C22").Select()
